# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old worker rows that are no longer part of this statement:
#  - old row 18 (JOSE LUIS MUÑOZ SEPULVEDA, will be re-added in row 17 below)
#  - old rows 19-24 (GUALBERTO EDUARDO VELASQUEZ TORDECILLA x6 periods)
# Deleting these 7 rows shifts the trailing "total" row (old row 25) up to
# row 18 (keeping its special border styling), and the signature lines
# (old rows 30-31) up to rows 23-24.
$ws.Rows("18:24").Delete()

# Update the header figures
$ws.Range("E11").Value = 115427
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 2

# Row 16: Alexander Camelo Consuegra
$ws.Range("C16").Value = "73182033"
$ws.Range("D16").Value = "ALEXANDER CAMELO CONSUEGRA"
$ws.Range("E16").Value = "2309"
$ws.Range("F16").Value = 1547
$ws.Range("G16").Value = 1423500

# Row 17: Jose Luis Muñoz Sepulveda
$ws.Range("C17").Value = "73162417"
$ws.Range("D17").Value = "JOSE LUIS MUÑOZ SEPULVEDA"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18 (total row, carried up from old row 25 - Jan Carlos Nieto Pedrozo):
# only the period value needs to be refreshed
$ws.Range("E18").Value = "2508"
